$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.44610000000003
$ws.Range("C5").Value = -13.8744
$ws.Range("E7").Value = 11.7933
$ws.Range("C9").Value = -11.6343
$ws.Range("C11").Value = -13.8089
$ws.Range("E11").Value = 13.3001
$ws.Range("A21").Value = -21.15500000000002
$ws.Range("C21").Value = -11.256
$ws.Range("E21").Value = 12.3651
$ws.Range("A23").Value = -21.36970000000003
$ws.Range("A25").Value = -22.51540000000003
